$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "122 Комп'ютерні науки"

# Update data rows 2-4 (header row stays the same except D2)
$ws.Range("D2").Value = "Країна 1"

$ws.Range("A3").Value = "Інформатика"
$ws.Range("B3").Value = "Комп'ютерних наук та кібернетики"
$ws.Range("C3").Value = "КНУ Шевченка"
$ws.Range("D3").Value = "Україна"
$ws.Range("E3").Value = "Прикладна математика"

$ws.Range("A4").Value = "rori"
$ws.Range("B4").Value = "Прикладна математика"
$ws.Range("C4").Value = "щ"
$ws.Range("D4").Value = "Англія"

# Update the active selection to E3
$ws.Range("E3").Select()
